$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Phase 1: clean up the original paragraph's run-fragmentation (remove the
# w:proofErr-induced "Shahida" run splits) by replacing each split span with
# itself; Word's Find/Replace naturally re-merges runs that end up with
# identical formatting and drops the now-pointless proofErr markers.
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("My client, Shahida Rahman, is an Author", $true, $false, $false, $false, $false, $true, 1, $false, "My client, Shahida Rahman, is an Author", 2) | Out-Null
$d.Content.Find.Execute("children. Shahida uses computers", $true, $false, $false, $false, $false, $true, 1, $false, "children. Shahida uses computers", 2) | Out-Null
$d.Content.Find.Execute("experience with computers. Shahida generally", $true, $false, $false, $false, $false, $true, 1, $false, "experience with computers. Shahida generally", 2) | Out-Null
$d.Content.Find.Execute("to Shahida, who sends the books off to print", $true, $false, $false, $false, $false, $true, 1, $false, "to Shahida, who sends the books off to print", 2) | Out-Null
$d.Content.Find.Execute("disorganised. Shahida would like", $true, $false, $false, $false, $false, $true, 1, $false, "disorganised. Shahida would like", 2) | Out-Null
$d.Content.Find.Execute("these. Shahida also wants", $true, $false, $false, $false, $false, $true, 1, $false, "these. Shahida also wants", 2) | Out-Null

# ---------------------------------------------------------------------------
# Phase 2: the "_GoBack" bookmark currently sits mid-sentence (between
# "...ready to print file" and "s to Shahida..."); it needs to move to the
# very end of the document once the new trailing paragraph exists. Drop it
# here and re-create it later in the correct spot.
# ---------------------------------------------------------------------------
$bm = $d.Bookmarks("_GoBack")
$bm.Delete()

# ---------------------------------------------------------------------------
# Phase 3: split the single paragraph into the six target paragraphs and add
# the new heading / body text.
# ---------------------------------------------------------------------------
$mainPara = $d.Paragraphs(1)

# --- New heading paragraph before everything: "1.1 Client Identification" ---
$startRange = $mainPara.Range.Duplicate
$startRange.Collapse(1)
$startRange.InsertParagraphBefore()
$heading1 = $d.Paragraphs(1)
$heading1.Range.Text = "1.1 Client Identification"
$heading1.Range.Font.Underline = 1

# $mainPara (the "My client, ..." text) is now Paragraphs(2).
$mainPara = $d.Paragraphs(2)

# --- New heading paragraph after the "My client..." paragraph: "1.2 Define the Current System" ---
$afterMain = $mainPara.Range.Duplicate
$afterMain.Collapse(0)
$afterMain.InsertParagraphAfter()
$heading2 = $d.Paragraphs(3)
$heading2.Range.Text = "1.2 Define the Current System"
$heading2.Range.Font.Underline = 1

# --- New body paragraph under 1.2 ---
$afterHeading2 = $heading2.Range.Duplicate
$afterHeading2.Collapse(0)
$afterHeading2.InsertParagraphAfter()
$body2 = $d.Paragraphs(4)
$body2.Range.Text = "The system that is currently being used consists of Shahida entering the book and its details into the spreadsheet. These details are taken from the enquiries that she receives via email, and include; author, book title, size, number of pages, hardback/paperback, mat or gloss, crème or white paper, font and font size. She also records their details in a separate spreadsheet, which includes their email, phone number, and address. Subsequently, Shahida waits for full payment and then sends the customer an invoice. She then contacts her editor and her illustrator to start work on the book. Shahida refers to her readily calculated prices for books in order to correctly price the book, in accordance to the book's details. Once the book is finished, the book is sent off to print, and the author receives 25 copies."

# --- New heading paragraph: "1.3 Problems" ---
$afterBody2 = $body2.Range.Duplicate
$afterBody2.Collapse(0)
$afterBody2.InsertParagraphAfter()
$heading3 = $d.Paragraphs(5)
$heading3.Range.Text = "1.3 Problems"
$heading3.Range.Font.Underline = 1

# --- New body paragraph under 1.3 ---
$afterHeading3 = $heading3.Range.Duplicate
$afterHeading3.Collapse(0)
$afterHeading3.InsertParagraphAfter()
$body3 = $d.Paragraphs(6)
$body3.Range.Text = "There are numerous problems with the current system. First of all, the usage of the spreadsheet makes it harder to find a customer and their details, and their book's details. This is because the spreadsheet is much disorganised. Furthermore, it is harder to keep track of the details of each book, meaning it is difficult to update the details of the book when necessary. Also, if the same author makes an enquiry about another book, her details must be entered into the spreadsheet again, which could cause inconsistencies in the data, because for instance, the customer may move house, meaning their address would need changing, and it would be difficult to find and update all entries where their address is recorded."

# ---------------------------------------------------------------------------
# Phase 4: re-create the "_GoBack" bookmark collapsed at the very end of the
# document. A collapsed bookmark whose position is the very last character
# offset of a paragraph confuses this host's Bookmarks.Add (it resets Start
# to 0), so append a throwaway character, anchor the bookmark just before it
# while that position is safely "mid-paragraph", then delete the throwaway
# character again.
# ---------------------------------------------------------------------------
$tail = $d.Content.Duplicate
$tail.Collapse(0)
$tail.InsertAfter("Z")

$endPos = $d.Content.End - 2
$bmRange = $d.Range($endPos, $endPos)
$d.Bookmarks.Add("_GoBack", $bmRange) | Out-Null

$placeholder = $d.Range($d.Content.End - 2, $d.Content.End - 1)
$placeholder.Delete()

Write-Output $d.Content.Text
